$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.459612070389937
$ws.Range("C2").Value = 10.29869402782916
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 21.22402822075207

$ws.Range("B3").Value = 0.3048080303191223
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 645.3272768299601
$ws.Range("G3").Value = 648.1053690801313
